$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "236.29") need to be
# forced to Text format first, otherwise Excel auto-converts them to numbers,
# which would not match the original inline-string (text) cell type.
$textCells = @(
    "D5",
    "D6",
    "D8",
    "D10",
    "D15",
    "D20",
    "D22",
    "D23",
    "D25",
    "D26",
    "D27",
    "D28",
    "D31",
    "D32",
    "D33",
    "D35",
    "D36",
    "D37",
    "D38",
    "D40",
    "D41",
    "D42",
    "D43",
    "D47",
    "D48",
    "D50",
    "D51",
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply all the updated values
$ws.Range("D2").Value = '94.200.64'
$ws.Range("E2").Value = '  +2.57%  '
$ws.Range("D3").Value = '3.073.64'
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '236.29'
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").Value = '609.70'
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("E7").Value = '  +1.91%  '
$ws.Range("D8").Value = '0.379'
$ws.Range("E8").Value = '  -2.07%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = '0.817'
$ws.Range("E10").Value = '  +11.07%  '
$ws.Range("D11").Value = '3.071.94'
$ws.Range("E11").Value = '  -0.55%  '
$ws.Range("E12").Value = '  -1.75%  '
$ws.Range("D13").Value = '93.996.49'
$ws.Range("E13").Value = '  +2.21%  '
$ws.Range("E14").Value = '  -2.24%  '
$ws.Range("D15").Value = '33.86'
$ws.Range("E15").Value = '  +0.05%  '
$ws.Range("E16").Value = '  -1.75%  '
$ws.Range("D17").Value = '3.643.96'
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").Value = '3.070.56'
$ws.Range("E18").Value = '  -0.68%  '
$ws.Range("E19").Value = '  -3.23%  '
$ws.Range("D20").Value = '14.36'
$ws.Range("E20").Value = '  -1.19%  '
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("D22").Value = '444.87'
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("D23").Value = '8.82'
$ws.Range("E23").Value = '  -4.61%  '
$ws.Range("E24").Value = '  -2.00%  '
$ws.Range("D25").Value = '8.36'
$ws.Range("E25").Value = '  +6.65%  '
$ws.Range("D26").Value = '5.51'
$ws.Range("E26").Value = '  -3.51%  '
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").Value = '11.96'
$ws.Range("E27").Value = '  +3.18%  '
$ws.Range("B28").Value = 'Litecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D28").Value = '84.61'
$ws.Range("E28").Value = '  -1.37%  '
$ws.Range("D29").Value = '3.237.31'
$ws.Range("E29").Value = '  -0.68%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").Value = '0.250'
$ws.Range("E31").Value = '  +7.56%  '
$ws.Range("D32").Value = '0.180'
$ws.Range("E32").Value = '  +7.48%  '
$ws.Range("D33").Value = '0.123'
$ws.Range("E33").Value = '  -6.05%  '
$ws.Range("E34").Value = '  +31.84%  '
$ws.Range("D35").Value = '8.96'
$ws.Range("E35").Value = '  -0.80%  '
$ws.Range("D36").Value = '7.56'
$ws.Range("E36").Value = '  -3.26%  '
$ws.Range("D37").Value = '25.45'
$ws.Range("E37").Value = '  -1.61%  '
$ws.Range("D38").Value = '0.151'
$ws.Range("E38").Value = '  -3.97%  '
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").Value = '479.19'
$ws.Range("E40").Value = '  -0.17%  '
$ws.Range("D41").Value = '24.05'
$ws.Range("E41").Value = '  +0.75%  '
$ws.Range("B42").Value = 'MantraDAO'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D42").Value = '3.79'
$ws.Range("E42").Value = '  -2.78%  '
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").Value = '0.436'
$ws.Range("E43").Value = '  +1.72%  '
$ws.Range("E44").Value = '  -2.44%  '
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("E46").Value = '  -5.32%  '
$ws.Range("D47").Value = '161.48'
$ws.Range("E47").Value = '  -1.85%  '
$ws.Range("D48").Value = '0.672'
$ws.Range("E48").Value = '  -1.79%  '
$ws.Range("E49").Value = '  -2.86%  '
$ws.Range("D50").Value = '43.65'
$ws.Range("B51").Value = 'FLOKI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D51").Value = '0.000269'
$ws.Range("E51").Value = '  +10.27%  '

# Restore the default (Normal) style on the text-forced cells so no stray
# cell-level number formatting is left behind
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
